$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Git Commit ID text (used in column AJ, "ScriptLatestRunVersion")
$oldCommitText = "IndicatorQuantiles.R, Git Commit ID: d77a77d64f72a744c78cd38270c72c5d9c8cd498"
$newCommitText = "IndicatorQuantiles.R, Git Commit ID: 24c5634628309d80791a95cb6332cf2c12927180"

$usedRange = $ws.UsedRange
$found = $usedRange.Find($oldCommitText)
if ($found -ne $null) {
    $firstAddress = $found.Address()
    do {
        $found.Value = $newCommitText
        $found = $usedRange.FindNext($found)
    } while ($found -ne $null -and $found.Address() -ne $firstAddress)
}

# Update the "pid" values in column AH from 11992 to 17548 for data rows 2-80
$lastRow = $ws.Cells.Item($ws.Rows.Count, 34).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 34)
    if ($cell.Value() -eq 11992) {
        $cell.Value = 17548
    }
}
